# Spreading out PPP over two quarters instead of three.
# Updates the "impact" (C), "total" (E) and "consumption" (H) columns
# for rows 87-94 (2021 Q1 through 2022 Q4) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (impact, total, consumption)
$updates = @{
    87 = @(5.547110201355259, 7.337069943294125, 5.190361128725651)
    88 = @(1.113172851774511, -3.445971813557383, -3.007484414320345)
    89 = @(-1.243179934994805, -5.818172229432157, -4.790084477699868)
    90 = @(-0.9876646317202245, -2.023584427185488, -1.887367382839073)
    91 = @(-3.379268200984162, -2.229344333761625, -1.622499104474916)
    92 = @(-3.034217630275573, -2.065769530723027, -1.996871803444562)
    93 = @(-1.751101248774911, -0.6857067034295107, -0.5787354101843698)
    94 = @(-1.290812433035287, -0.182429164226992, 0.009996642839077443)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]   # column C - impact
    $ws.Cells.Item($row, 5).Value = $vals[1]   # column E - total
    $ws.Cells.Item($row, 8).Value = $vals[2]   # column H - consumption
}
